# Apply cryptos list update (coinranking prices/volumes refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "585.90"); force
# a Text number format first so COM does not coerce/trim it into a float.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.435.20"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.609.77"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.90"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.88"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.604.26"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.632"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.665"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.17"
$ws.Range("E12").Value = "  -4.24%  "
$ws.Range("E13").Value = "  +8.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.75"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.186.64"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.605.72"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.379.78"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.47"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "20.02"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.93"
$ws.Range("E24").Value = "  -7.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.40"
$ws.Range("E25").Value = "  +7.06%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.11"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.45"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "579.68"
$ws.Range("E35").Value = "  -9.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.05"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0819"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.30"
$ws.Range("E40").Value = "  +21.12%  "
$ws.Range("E41").Value = "  +6.13%  "
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  -6.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.226.60"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.09"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0448"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("E47").Value = "  +5.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.139"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.06%  "
